$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '55.329.48'
$ws.Range("E2").Value = "'" + '  +8.02%  '

$ws.Range("D3").Value = "'" + '2.436.38'
$ws.Range("E3").Value = "'" + '  +8.20%  '

$ws.Range("D4").Value = "'" + '0.999'
$ws.Range("E4").Value = "'" + '  +0.07%  '

$ws.Range("D5").Value = "'" + '480.44'
$ws.Range("E5").Value = "'" + '  +11.80%  '

$ws.Range("D6").Value = "'" + '140.14'
$ws.Range("E6").Value = "'" + '  +22.36%  '

$ws.Range("E7").Value = "'" + '  +0.10%  '

$ws.Range("D8").Value = "'" + '0.505'
$ws.Range("E8").Value = "'" + '  +12.66%  '

$ws.Range("D9").Value = "'" + '2.457.92'
$ws.Range("E9").Value = "'" + '  +8.99%  '

$ws.Range("D10").Value = "'" + '0.0964'
$ws.Range("E10").Value = "'" + '  +16.23%  '

$ws.Range("D11").Value = "'" + '5.48'
$ws.Range("E11").Value = "'" + '  +8.34%  '

$ws.Range("D12").Value = "'" + '0.325'
$ws.Range("E12").Value = "'" + '  +11.71%  '

$ws.Range("D13").Value = "'" + '0.123'
$ws.Range("E13").Value = "'" + '  +3.02%  '

$ws.Range("D14").Value = "'" + '2.854.85'
$ws.Range("E14").Value = "'" + '  +9.11%  '

$ws.Range("D15").Value = "'" + '55.289.70'
$ws.Range("E15").Value = "'" + '  +8.07%  '

$ws.Range("E16").Value = "'" + '  +13.69%  '

$ws.Range("D17").Value = "'" + '0.0000135'
$ws.Range("E17").Value = "'" + '  +21.00%  '

$ws.Range("D18").Value = "'" + '2.460.01'
$ws.Range("E18").Value = "'" + '  +9.65%  '

$ws.Range("D19").Value = "'" + '4.36'
$ws.Range("E19").Value = "'" + '  +14.11%  '

$ws.Range("E20").Value = "'" + '  +20.12%  '

$ws.Range("D21").Value = "'" + '314.69'
$ws.Range("E21").Value = "'" + '  +9.68%  '

$ws.Range("E22").Value = "'" + '  +0.55%  '

$ws.Range("E23").Value = "'" + '  +16.55%  '

$ws.Range("D24").Value = "'" + '57.44'
$ws.Range("E24").Value = "'" + '  +10.48%  '

$ws.Range("E25").Value = "'" + '  +0.41%  '

$ws.Range("E26").Value = "'" + '  +13.08%  '

$ws.Range("E27").Value = "'" + '  +21.96%  '

$ws.Range("D28").Value = "'" + '2.540.08'
$ws.Range("E28").Value = "'" + '  +9.24%  '

$ws.Range("D29").Value = "'" + '7.37'
$ws.Range("E29").Value = "'" + '  +13.19%  '

$ws.Range("D30").Value = "'" + '0.0₃0770'
$ws.Range("E30").Value = "'" + '  +26.65%  '

$ws.Range("E31").Value = "'" + '  +0.07%  '

$ws.Range("D32").Value = "'" + '148.51'
$ws.Range("E32").Value = "'" + '  +4.67%  '

$ws.Range("D33").Value = "'" + '17.92'
$ws.Range("E33").Value = "'" + '  +11.05%  '

$ws.Range("E34").Value = "'" + '  +15.06%  '

$ws.Range("E35").Value = "'" + '  +14.86%  '

$ws.Range("E36").Value = "'" + '  +19.58%  '

$ws.Range("E37").Value = "'" + '  +11.32%  '

$ws.Range("D38").Value = "'" + '0.841'
$ws.Range("E38").Value = "'" + '  +17.09%  '

$ws.Range("D39").Value = "'" + '33.49'
$ws.Range("E39").Value = "'" + '  +6.80%  '

$ws.Range("D40").Value = "'" + '0.994'
$ws.Range("E40").Value = "'" + '  +0.12%  '

$ws.Range("E41").Value = "'" + '  +10.61%  '

$ws.Range("D42").Value = "'" + '3.41'
$ws.Range("E42").Value = "'" + '  +13.79%  '

$ws.Range("E43").Value = "'" + '  +13.93%  '

$ws.Range("D44").Value = "'" + '1.29'
$ws.Range("E44").Value = "'" + '  +19.46%  '

$ws.Range("D45").Value = "'" + '10.12'
$ws.Range("E45").Value = "'" + '  -0.02%  '

$ws.Range("B46").Value = "'" + 'RenderToken'
$ws.Range("C46").Value = "'" + 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = "'" + '4.63'
$ws.Range("E46").Value = "'" + '  +22.82%  '

$ws.Range("B47").Value = "'" + 'Bittensor'
$ws.Range("C47").Value = "'" + 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").Value = "'" + '253.58'
$ws.Range("E47").Value = "'" + '  +35.54%  '

$ws.Range("D48").Value = "'" + '0.0896'
$ws.Range("E48").Value = "'" + '  +15.45%  '

$ws.Range("E49").Value = "'" + '  +14.09%  '

$ws.Range("D50").Value = "'" + '1.919.20'
$ws.Range("E50").Value = "'" + '  +5.03%  '

$ws.Range("D51").Value = "'" + '17.17'
$ws.Range("E51").Value = "'" + '  +14.36%  '
